$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.370.44'
$ws.Range('E2').Value = '  +2.39%  '
$ws.Range('D3').Value = '3.345.17'
$ws.Range('E3').Value = '  +3.38%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'192.88"
$ws.Range('E5').Value = '  +4.98%  '
$ws.Range('D6').Value = "'591.23"
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.91%  '
$ws.Range('E9').Value = '  +3.27%  '
$ws.Range('D10').Value = "'6.75"
$ws.Range('E10').Value = '  +2.40%  '
$ws.Range('D11').Value = "'0.422"
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('D12').Value = '3.928.85'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = "'28.19"
$ws.Range('E14').Value = '  +2.24%  '
$ws.Range('D15').Value = '69.374.41'
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('E16').Value = '  +1.68%  '
$ws.Range('D17').Value = '3.370.78'
$ws.Range('E17').Value = '  +4.95%  '
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('E19').Value = '  +2.43%  '
$ws.Range('D20').Value = "'429.32"
$ws.Range('E20').Value = '  +8.42%  '
$ws.Range('D22').Value = "'73.15"
$ws.Range('E22').Value = '  +2.94%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('E25').Value = '  +3.15%  '
$ws.Range('E26').Value = '  +3.51%  '
$ws.Range('D27').Value = "'9.60"
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('D28').Value = "'1.02"
$ws.Range('E28').Value = '  +2.27%  '
$ws.Range('E29').Value = '  +2.26%  '
$ws.Range('D30').Value = "'23.03"
$ws.Range('E30').Value = '  +1.80%  '
$ws.Range('D31').Value = "'5.59"
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('E32').Value = '  +1.96%  '
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('D34').Value = "'0.999"
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').Value = "'164.84"
$ws.Range('E35').Value = '  +1.92%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = "'1.52"
$ws.Range('E36').Value = '  +2.89%  '
$ws.Range('E37').Value = '  +1.95%  '
$ws.Range('D38').Value = "'27.01"
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = "'6.50"
$ws.Range('E41').Value = '  +0.85%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.754.84'
$ws.Range('E42').Value = '  +5.77%  '
$ws.Range('E43').Value = '  +2.24%  '
$ws.Range('D44').Value = "'41.20"
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = "'0.0686"
$ws.Range('E45').Value = '  +0.64%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = "'344.53"
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = "'25.29"
$ws.Range('E47').Value = '  +1.78%  '
$ws.Range('E48').Value = '  +2.00%  '
$ws.Range('D49').Value = "'32.64"
$ws.Range('E49').Value = '  +5.83%  '
$ws.Range('E50').Value = '  +3.79%  '
$ws.Range('E51').Value = '  +0.23%  '
